$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the value of E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Set the selection to E8 (as recorded by Excel when the file was last saved)
$ws.Range("E8").Select()
